# Apply the "Updated remaining queries for C3DC" edit:
#  - Update the JOIN conditions in every SQL query cell (columns B and C on
#    rows 2-7 of the active sheet) so that the generic ".id" columns are
#    replaced with the fully qualified "<table>_id" columns.
#  - Resize column C to a fixed width of 69 (removing the "best fit"
#    autofit sizing that was previously applied).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells containing the SQL text that needs the JOIN-column rename.
$cells = @("C2", "B2", "B3", "B4", "B5", "B6", "B7")

foreach ($cellAddr in $cells) {
    $rng = $ws.Range($cellAddr)
    $text = $rng.Value()

    $text = $text.Replace('std.id = prt."study.id"', 'std.study_id = prt."study.study_id"')
    $text = $text.Replace('prt.id = dgn."participant.id"', 'prt.participant_id = dgn."participant.participant_id"')
    $text = $text.Replace('prt.id = trt."participant.id"', 'prt.participant_id = trt."participant.participant_id"')
    $text = $text.Replace('prt.id = trr."participant.id"', 'prt.participant_id = trr."participant.participant_id"')
    $text = $text.Replace('prt.id = srv."participant.id"', 'prt.participant_id = srv."participant.participant_id"')
    $text = $text.Replace('std.id = rfs."study.id"', 'std.study_id = rfs."study.study_id"')

    $rng.Value = $text
}

# Column C: change from "best fit" auto-sized width to a fixed width of 69.
$ws.Columns.Item(3).ColumnWidth = 68.15
